$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1480.0922
$ws.Range("I15").Value = 1480.0922
$ws.Range("K15").Value = 4440.2766
$ws.Range("M15").Value = -4271.2766
$ws.Range("H46").Value = 3330.3333
$ws.Range("I46").Value = 2995.5
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 8986.5
$ws.Range("L46").Value = 12000
$ws.Range("M46").Value = -8867.5
$ws.Range("N46").Value = -12238
$ws.Range("H60").Value = 3330.3333
$ws.Range("I60").Value = 2995.5
$ws.Range("J60").Value = 4000
$ws.Range("K60").Value = 8986.5
$ws.Range("L60").Value = 12000
$ws.Range("M60").Value = -8502.5
$ws.Range("N60").Value = -12968
$ws.Range("H106").Value = 2217.9
$ws.Range("I106").Value = 2023.1875
$ws.Range("K106").Value = 2023.1875
$ws.Range("M106").Value = -1392.1875
$ws.Range("H116").Value = 5199.4546
$ws.Range("I116").Value = 5205.6
$ws.Range("K116").Value = 5205.6
$ws.Range("M116").Value = -1763.6
$ws.Range("H132").Value = 7558051.5
$ws.Range("I132").Value = 7977893
$ws.Range("J132").Value = 900
$ws.Range("K132").Value = 23933679
$ws.Range("L132").Value = 2700
$ws.Range("M132").Value = -23931149
$ws.Range("N132").Value = -7760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24774.666
$ws.Range("J32").Value = 10000
$ws.Range("L32").Value = 10000
$ws.Range("N32").Value = -10574
$ws.Range("H37").Value = 19374.5
$ws.Range("J37").Value = 26666
$ws.Range("L37").Value = 26666
$ws.Range("N37").Value = -27212
$ws.Range("H45").Value = 3587.2666
$ws.Range("I45").Value = 2346.5715
$ws.Range("J45").Value = 4672.875
$ws.Range("K45").Value = 2346.5715
$ws.Range("L45").Value = 4672.875
$ws.Range("M45").Value = -1969.5715
$ws.Range("N45").Value = -5426.875
$ws.Range("H61").Value = 3684.6562
$ws.Range("I61").Value = 1143.6666
$ws.Range("J61").Value = 8535.637000000001
$ws.Range("K61").Value = 1143.6666
$ws.Range("L61").Value = 8535.637000000001
$ws.Range("M61").Value = -931.6666
$ws.Range("N61").Value = -8959.637000000001
$ws.Range("H74").Value = 1004585
$ws.Range("I74").Value = 2003003.6
$ws.Range("K74").Value = 2003003.6
$ws.Range("M74").Value = -2002129.6
$ws.Range("H77").Value = 1004585
$ws.Range("I77").Value = 2003003.6
$ws.Range("K77").Value = 10015018
$ws.Range("M77").Value = -10010650
$ws.Range("H86").Value = 50000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = $null
$ws.Range("H89").Value = 50000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = $null
$ws.Range("H122").Value = 3502
$ws.Range("I122").Value = 3737.3333
$ws.Range("K122").Value = 11211.9999
$ws.Range("M122").Value = -8761.999899999999
$ws.Range("H132").Value = 2220.6
$ws.Range("I132").Value = 1479.2222
$ws.Range("J132").Value = 3332.6667
$ws.Range("K132").Value = 4437.6666
$ws.Range("L132").Value = 9998.000100000001
$ws.Range("M132").Value = -1907.6666
$ws.Range("N132").Value = -15058.0001
$ws.Range("H134").Value = 63529.75
$ws.Range("J134").Value = 63529.75
$ws.Range("L134").Value = 63529.75
$ws.Range("N134").Value = -73669.75
$ws.Range("H135").Value = 35571.2
$ws.Range("J135").Value = 35571.2
$ws.Range("L135").Value = 35571.2
$ws.Range("N135").Value = -45711.2
$ws.Range("H136").Value = 3684.6562
$ws.Range("I136").Value = 1143.6666
$ws.Range("J136").Value = 8535.637000000001
$ws.Range("K136").Value = 3430.9998
$ws.Range("L136").Value = 25606.911
$ws.Range("M136").Value = -880.9998000000001
$ws.Range("N136").Value = -30706.911

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 27268.678
$ws.Range("J20").Value = 1469.3334
$ws.Range("L20").Value = 1469.3334
$ws.Range("N20").Value = -1963.3334
$ws.Range("H22").Value = 379
$ws.Range("I22").Value = 379
$ws.Range("K22").Value = 379
$ws.Range("M22").Value = -206
$ws.Range("H123").Value = 16086.667
$ws.Range("H134").Value = 3023.2334
$ws.Range("I134").Value = 2720.077
$ws.Range("K134").Value = 8160.231000000001
$ws.Range("M134").Value = -5625.231000000001
$ws.Range("H135").Value = 90468
$ws.Range("J135").Value = 90468
$ws.Range("L135").Value = 90468
$ws.Range("N135").Value = -100608

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 609
$ws.Range("I22").Value = 636.75
$ws.Range("K22").Value = 636.75
$ws.Range("M22").Value = -286.75
$ws.Range("H62").Value = 8387.200000000001
$ws.Range("I62").Value = 3388.6
$ws.Range("J62").Value = 13385.8
$ws.Range("K62").Value = 3388.6
$ws.Range("L62").Value = 13385.8
$ws.Range("M62").Value = -2764.6
$ws.Range("N62").Value = -14633.8
$ws.Range("H65").Value = 8387.200000000001
$ws.Range("I65").Value = 3388.6
$ws.Range("J65").Value = 13385.8
$ws.Range("K65").Value = 16943
$ws.Range("L65").Value = 66929
$ws.Range("M65").Value = -13823
$ws.Range("N65").Value = -73169
$ws.Range("H86").Value = 84778.664
$ws.Range("J86").Value = 26626
$ws.Range("L86").Value = 26626
$ws.Range("N86").Value = -28872
$ws.Range("H89").Value = 84778.664
$ws.Range("J89").Value = 26626
$ws.Range("L89").Value = 133130
$ws.Range("N89").Value = -144362
$ws.Range("H99").Value = 12878.5
$ws.Range("I99").Value = 15995
$ws.Range("J99").Value = 11839.667
$ws.Range("K99").Value = 15995
$ws.Range("L99").Value = 11839.667
$ws.Range("M99").Value = -14497
$ws.Range("N99").Value = -14835.667
$ws.Range("H126").Value = 12878.5
$ws.Range("I126").Value = 15995
$ws.Range("J126").Value = 11839.667
$ws.Range("K126").Value = 47985
$ws.Range("L126").Value = 35519.001
$ws.Range("M126").Value = -45515
$ws.Range("N126").Value = -40459.001
$ws.Range("H132").Value = 38622.312
$ws.Range("I132").Value = 41111.133
$ws.Range("K132").Value = 123333.399
$ws.Range("M132").Value = -120803.399
$ws.Range("H134").Value = 3677.3333
$ws.Range("I134").Value = 3637
$ws.Range("K134").Value = 10911
$ws.Range("M134").Value = -8376

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 323.85715
$ws.Range("I8").Value = 323.85715
$ws.Range("K8").Value = 971.5714499999999
$ws.Range("M8").Value = -832.5714499999999
$ws.Range("H45").Value = 3749.1667
$ws.Range("J45").Value = 3693
$ws.Range("L45").Value = 11079
$ws.Range("N45").Value = -12143
$ws.Range("H70").Value = 4869.478
$ws.Range("H73").Value = 4869.478
$ws.Range("H138").Value = 4384.9287
$ws.Range("I138").Value = 1336.4546
$ws.Range("J138").Value = 15562.667
$ws.Range("K138").Value = 4009.3638
$ws.Range("L138").Value = 46688.001
$ws.Range("M138").Value = 1130.6362
$ws.Range("N138").Value = -56968.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 116.8
$ws.Range("I2").Value = 193.75
$ws.Range("J2").Value = 80.588234
$ws.Range("K2").Value = 193.75
$ws.Range("L2").Value = 80.588234
$ws.Range("M2").Value = -80.75
$ws.Range("N2").Value = -306.588234
$ws.Range("H113").Value = 1776.5555
$ws.Range("I113").Value = 818.8
$ws.Range("K113").Value = 818.8
$ws.Range("M113").Value = 1351.2
$ws.Range("H122").Value = 3292.6428
$ws.Range("I122").Value = 3235.0527
$ws.Range("K122").Value = 9705.158100000001
$ws.Range("M122").Value = -7255.158100000001
$ws.Range("H132").Value = 4066.0715
$ws.Range("I132").Value = 4393.5
$ws.Range("J132").Value = 3247.5
$ws.Range("K132").Value = 13180.5
$ws.Range("L132").Value = 9742.5
$ws.Range("M132").Value = -10650.5
$ws.Range("N132").Value = -14802.5
$ws.Range("H135").Value = 129998.75
$ws.Range("J135").Value = 129998.75
$ws.Range("L135").Value = 129998.75
$ws.Range("N135").Value = -140138.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2993.261
$ws.Range("I40").Value = 2881.4211
$ws.Range("K40").Value = 2881.4211
$ws.Range("M40").Value = -2745.4211
$ws.Range("H61").Value = 522.1
$ws.Range("I61").Value = 531.25
$ws.Range("J61").Value = 508.375
$ws.Range("K61").Value = 531.25
$ws.Range("L61").Value = 508.375
$ws.Range("M61").Value = -329.25
$ws.Range("N61").Value = -912.375
$ws.Range("H113").Value = 522.1
$ws.Range("I113").Value = 531.25
$ws.Range("J113").Value = 508.375
$ws.Range("K113").Value = 531.25
$ws.Range("L113").Value = 508.375
$ws.Range("M113").Value = 1638.75
$ws.Range("N113").Value = -4848.375
$ws.Range("H122").Value = 4747.6665
$ws.Range("I122").Value = 3585.6667
$ws.Range("K122").Value = 10757.0001
$ws.Range("M122").Value = -8307.000100000001
$ws.Range("H132").Value = 5970.8887
$ws.Range("I132").Value = 5734
$ws.Range("J132").Value = 6533.5
$ws.Range("K132").Value = 17202
$ws.Range("L132").Value = 19600.5
$ws.Range("M132").Value = -14672
$ws.Range("N132").Value = -24660.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3412.6453
$ws.Range("I126").Value = 3007.4
$ws.Range("K126").Value = 9022.200000000001
$ws.Range("M126").Value = -6552.200000000001
$ws.Range("H136").Value = 11371.392
$ws.Range("I136").Value = 13353.569
$ws.Range("K136").Value = 40060.70699999999
$ws.Range("M136").Value = -37510.70699999999
